# Fruta / hortaliza, semanal
# Update the weekly price rows (dates/volumes/prices/origin) for the
# "Vega Central Mapocho de Santiago - Higo" sheet so that each data row
# (4-9, 12-18) reflects its updated weekly figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 44299
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range('R4').Value = 'Provincia de Santiago'
$ws.Range("S4").Value = 2143
$ws.Range("D5").Value = 44299
$ws.Range("M5").Value = 75
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range('R5').Value = 'Provincia de Santiago'
$ws.Range("S5").Value = 1714
$ws.Range("D6").Value = 44302
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 15000
$ws.Range("S6").Value = 2143
$ws.Range("D7").Value = 44302
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 12000
$ws.Range("S7").Value = 1714
$ws.Range("D8").Value = 44980
$ws.Range("N8").Value = 16000
$ws.Range("O8").Value = 16000
$ws.Range("P8").Value = 16000
$ws.Range('R8').Value = 'Región Metropolitana'
$ws.Range("S8").Value = 2286
$ws.Range("D9").Value = 44980
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 13000
$ws.Range("O9").Value = 13000
$ws.Range("P9").Value = 13000
$ws.Range('R9').Value = 'Región Metropolitana'
$ws.Range("S9").Value = 1857
$ws.Range("D12").Value = 44300
$ws.Range("M12").Value = 100
$ws.Range("D13").Value = 44300
$ws.Range("M13").Value = 80
$ws.Range("D14").Value = 44301
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 14000
$ws.Range("O14").Value = 14000
$ws.Range("P14").Value = 14000
$ws.Range('Q14').Value = '$/bandeja 7 kilos'
$ws.Range("S14").Value = 2000
$ws.Range("T14").Value = 7
$ws.Range("D15").Value = 44301
$ws.Range('L15').Value = 'Segunda'
$ws.Range("M15").Value = 80
$ws.Range("N15").Value = 12000
$ws.Range("O15").Value = 12000
$ws.Range("P15").Value = 12000
$ws.Range("S15").Value = 1714
$ws.Range("D16").Value = 44322
$ws.Range('L16').Value = 'Primera'
$ws.Range("M16").Value = 45
$ws.Range("D17").Value = 44322
$ws.Range('L17').Value = 'Segunda'
$ws.Range("M17").Value = 80
$ws.Range("N17").Value = 8000
$ws.Range("O17").Value = 8000
$ws.Range("P17").Value = 8000
$ws.Range("S17").Value = 1143
$ws.Range("D18").Value = 44971
$ws.Range('L18').Value = 'Primera'
$ws.Range("M18").Value = 25
$ws.Range("N18").Value = 15000
$ws.Range("O18").Value = 15000
$ws.Range("P18").Value = 15000
$ws.Range('Q18').Value = '$/bandeja 5 kilos'
$ws.Range("S18").Value = 3000
$ws.Range("T18").Value = 5
